# Update the dSF column (column F) values on the active worksheet.
# Each pair is (row number, new value) corresponding to the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -6
    4  = -6
    5  = -7
    7  = -1
    11 = 2
    17 = -1
    19 = -2
    21 = -2
    22 = 0
    23 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
